# Apply the changes described in the commit "Cambios de los excel e icono de usuarios"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2 / B3: "Negocios" -> "Sistemas"
$ws.Range("B2").Value = "Sistemas"
$ws.Range("B3").Value = "Sistemas"

# Update the active selection to match the saved view (was J5, now E7)
$ws.Range("E7").Select()
